$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1653195549
$ws.Range("G3").Value = 1653195544
$ws.Range("G5").Value = 1653195556
$ws.Range("G6").Value = 1653176010
$ws.Range("G7").Value = 1653186311
